$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 34, shifting existing rows 34:84 down to 35:85
$ws.Rows.Item(34).Insert()

# Helper: write a value as plain TEXT (no number/date auto-detection, no
# extra NumberFormat/quote-prefix style) by building it as a text formula in
# a scratch cell, then Paste-Special-Values-only into the destination - this
# mirrors how the sheet's existing inline-string numeric-looking cells
# (case ids, OT ids, comuna numbers, etc.) are stored.
function Set-TextValue($addr, $text) {
    $scratch = $ws.Range("Z1")
    $escaped = $text.Replace("""", """""")
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-TextValue "A34" "6045"
Set-TextValue "B34" "2/7/2025"
$ws.Range("C34").Value = "GALLARDO, ANGEL AV. 213"
Set-TextValue "D34" "15"
Set-TextValue "E34" "803607430"
$ws.Range("F34").Value = "PEBCOM"
$ws.Range("G34").Value = "Pendiente"
$ws.Range("H34").Value = "Era el caso 4852 volvio a entrar por estar mal cementada la base volver a reparar"
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = "Cambio"
$ws.Range("K34").Value = "Sin equipos"
$ws.Range("L34").Value = "Pasante"
$ws.Range("M34").Value = -58.435452
$ws.Range("N34").Value = -34.603627
$ws.Range("O34").Value = "Fuera de operaciones"
$ws.Range("P34").Value = "No clasificado, consultar con mantenimiento"
